# Applies the commit's edit:
#   - Appends a closing sentence to the paragraph that ends in
#     "...changing nothing."
#   - Removes the three trailing paragraphs that discuss the Room Service
#     bug ("We ran into problems...", the "Basically the Value Guest..."
#     bullet, and "The problem lies in...").
#   - Keeps the _GoBack bookmark alive, relocating it to sit right before
#     the paragraph mark of the (now final) paragraph, exactly where it
#     ends up in the target document.

$d = $word.ActiveDocument

$newSentence = " With the downloading of information from the application and uploading from the database fully operational, this system is complete."

# --- Step 1: append the new sentence to the end of the paragraph that ---
# --- currently ends in "...changing nothing."                          ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*changing nothing.*") {
        $target = $candidate
        break
    }
}

$insertPos = $target.Range.End - 1
$insertPoint = $d.Range($insertPos, $insertPos)
$insertPoint.InsertAfter($newSentence)

# --- Step 2: delete the following three paragraphs entirely (the ---
# --- "We ran into problems...", the bulleted "Basically..." item, ---
# --- and "The problem lies..." paragraphs).                       ---
$firstToRemove = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*We ran into problems with the insert/update*") {
        $firstToRemove = $candidate
        break
    }
}

$lastToRemove = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*The problem lies in our insert/update*") {
        $lastToRemove = $candidate
    }
}

$removeRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$removeRange.Delete()

# --- Step 3: recreate the _GoBack bookmark, collapsed, right before ---
# --- the paragraph mark of the paragraph we just extended (which is ---
# --- now the final paragraph in the body).                          ---
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$markPos = $finalPara.Range.End - 1

# Placing a collapsed bookmark directly at a position immediately before a
# paragraph mark is unreliable, so two throw-away placeholder characters
# are inserted, the bookmark is anchored safely between them, and then the
# placeholders are deleted one at a time so the bookmark naturally settles
# into the correct collapsed position right before the paragraph mark.
$placeholder = $d.Range($markPos, $markPos)
$placeholder.InsertAfter("XY")

$midPos = $markPos + 1
$bookmarkRange = $d.Range($midPos, $midPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$afterChar = $d.Range($midPos, $midPos + 1)
$afterChar.Delete()

$beforeChar = $d.Range($midPos - 1, $midPos)
$beforeChar.Delete()
